$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously empty Accuracy/FAR/FRR cells for row 4 with text
# values (enter via formula so Excel keeps them as text, then convert the
# formulas to static values so the results are stored as plain text cells)
$ws.Range("B4").Formula = "=""0.8027"""
$ws.Range("C4").Formula = "=""0.2953"""
$ws.Range("D4").Formula = "=""0.0990"""
$ws.Range("B4:D4").Copy()
$ws.Range("B4:D4").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

# Copy the styling used in row 3 (B3:D3) onto the new cells (B4:D4)
$ws.Range("B3:D3").Copy()
$ws.Range("B4:D4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Select E7 to mirror the final cursor position
$ws.Range("E7").Select()
